$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("education")

# Correct the university name: "Universidad Jaime I " -> "Universitat Jaume I"
$ws.Range("C5").Value = "Universitat Jaume I"

# Update the active selection as recorded in the saved sheet view
$ws.Activate()
$ws.Range("C11").Select()
